$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45-151 down to 46-152.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new record's data.
$ws.Cells.Item(45, 1).Value = 7
$ws.Cells.Item(45, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(45, 3).Value = "Ñuble"
$ws.Cells.Item(45, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(45, 5).Value = 16
$ws.Cells.Item(45, 6).Value = 100112017
$ws.Cells.Item(45, 7).Value = "Apio"
$ws.Cells.Item(45, 8).Value = "Americana (o)"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 60
$ws.Cells.Item(45, 11).Value = 8000
$ws.Cells.Item(45, 12).Value = 8500
$ws.Cells.Item(45, 13).Value = 8250
$ws.Cells.Item(45, 14).Value = "`$/docena de matas"
$ws.Cells.Item(45, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(45, 16).Value = 1375
$ws.Cells.Item(45, 17).Value = 6
$ws.Cells.Item(45, 18).Value = "Hortaliza"
